$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string that must stay TEXT
# (matches original inlineStr formatting such as "0.999", "6.62", etc.)
$textForceCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D14", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D33", "D37", "D38", "D39", "D41", "D43", "D46", "D48", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated Price (D) values
$ws.Range("D5").Value = "580.21"
$ws.Range("D6").Value = "168.45"
$ws.Range("D7").Value = "0.999"
$ws.Range("D10").Value = "6.62"
$ws.Range("D11").Value = "0.154"
$ws.Range("D12").Value = "0.482"
$ws.Range("D14").Value = "36.49"
$ws.Range("D20").Value = "16.13"
$ws.Range("D21").Value = "465.46"
$ws.Range("D23").Value = "7.53"
$ws.Range("D24").Value = "83.25"
$ws.Range("D25").Value = "2.36"
$ws.Range("D26").Value = "12.90"
$ws.Range("D27").Value = "10.14"
$ws.Range("D29").Value = "7.97"
$ws.Range("D33").Value = "28.14"
$ws.Range("D37").Value = "5.89"
$ws.Range("D38").Value = "2.14"
$ws.Range("D39").Value = "46.72"
$ws.Range("D41").Value = "50.20"
$ws.Range("D43").Value = "8.69"
$ws.Range("D46").Value = "383.50"
$ws.Range("D48").Value = "134.63"
$ws.Range("D50").Value = "24.66"

# Restore default (Normal) style so no stray number-format style id is left on the cell
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining Price (D) / Volume(1h) (E) updates (values that are naturally text)
$ws.Range("D2").Value = "66.997.71"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "3.087.09"
$ws.Range("E3").Value = "  +4.63%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("E6").Value = "  +5.74%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.082.63"
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("E12").Value = "  +5.11%  "
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("E14").Value = "  +6.54%  "
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "3.594.93"
$ws.Range("E16").Value = "  +4.45%  "
$ws.Range("D17").Value = "66.903.74"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "3.083.16"
$ws.Range("E19").Value = "  +4.38%  "
$ws.Range("E20").Value = "  +7.80%  "
$ws.Range("E21").Value = "  +4.40%  "
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("E23").Value = "  +4.03%  "
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +6.87%  "
$ws.Range("E26").Value = "  +6.74%  "
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("E38").Value = "  +7.80%  "
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("E40").Value = "  +6.97%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").Value = "2.763.60"
$ws.Range("E47").Value = "  +1.98%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("E50").Value = "  +6.45%  "
$ws.Range("E51").Value = "  +2.57%  "
